# Update the "Förändrad" (column C) date values from 2023-09-01 (45170)
# to 2023-09-05 (45174) for all data rows (rows 2 through 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
